$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 10870394
$ws.Range("I33").Value = 13889330
$ws.Range("K33").Value = 13889330
$ws.Range("M33").Value = -13889101
$ws.Range("H98").Value = 1463.4572
$ws.Range("I98").Value = 1277.9062
$ws.Range("K98").Value = 1277.9062
$ws.Range("M98").Value = 220.0938000000001
$ws.Range("H122").Value = 1463.4572
$ws.Range("I122").Value = 1277.9062
$ws.Range("K122").Value = 3833.7186
$ws.Range("M122").Value = -1383.7186

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 32307.416
$ws.Range("I61").Value = 28415.5
$ws.Range("K61").Value = 28415.5
$ws.Range("M61").Value = -28203.5
$ws.Range("H102").Value = 8417.691999999999
$ws.Range("I102").Value = 754.4400000000001
$ws.Range("K102").Value = 754.4400000000001
$ws.Range("M102").Value = 867.5599999999999
$ws.Range("H132").Value = 2641440.5
$ws.Range("I132").Value = 4530.85
$ws.Range("J132").Value = 5571340.5
$ws.Range("K132").Value = 13592.55
$ws.Range("L132").Value = 16714021.5
$ws.Range("M132").Value = -11062.55
$ws.Range("N132").Value = -16719081.5
$ws.Range("H136").Value = 32307.416
$ws.Range("I136").Value = 28415.5
$ws.Range("K136").Value = 85246.5
$ws.Range("M136").Value = -82696.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 47477.332
$ws.Range("I20").Value = 40165.07
$ws.Range("K20").Value = 40165.07
$ws.Range("M20").Value = -39918.07
$ws.Range("H31").Value = 212.5
$ws.Range("I31").Value = 212.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 212.5
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("M31").Value = 39.5
$ws.Range("H86").Value = 1992.32
$ws.Range("I86").Value = 1657.3422
$ws.Range("J86").Value = 3053.0833
$ws.Range("K86").Value = 1657.3422
$ws.Range("L86").Value = 3053.0833
$ws.Range("M86").Value = -534.3422
$ws.Range("N86").Value = -5299.0833
$ws.Range("H89").Value = 1992.32
$ws.Range("I89").Value = 1657.3422
$ws.Range("J89").Value = 3053.0833
$ws.Range("K89").Value = 8286.710999999999
$ws.Range("L89").Value = 15265.4165
$ws.Range("M89").Value = -2670.710999999999
$ws.Range("N89").Value = -26497.4165
$ws.Range("H99").Value = 1049.75
$ws.Range("I99").Value = 999.7273
$ws.Range("J99").Value = 1600
$ws.Range("K99").Value = 999.7273
$ws.Range("L99").Value = 1600
$ws.Range("M99").Value = 498.2727
$ws.Range("N99").Value = -4596
$ws.Range("H134").Value = 24577.072
$ws.Range("I134").Value = 13694
$ws.Range("J134").Value = 39087.832
$ws.Range("K134").Value = 41082
$ws.Range("L134").Value = 117263.496
$ws.Range("M134").Value = -38547
$ws.Range("N134").Value = -122333.496

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 211.42105
$ws.Range("J7").Value = 255.28572
$ws.Range("L7").Value = 255.28572
$ws.Range("N7").Value = -481.28572
$ws.Range("H58").Value = 11620.02
$ws.Range("I58").Value = 6825.0835
$ws.Range("J58").Value = 13095.385
$ws.Range("K58").Value = 6825.0835
$ws.Range("L58").Value = 13095.385
$ws.Range("M58").Value = -6622.0835
$ws.Range("N58").Value = -13501.385
$ws.Range("H86").Value = 3426.9092
$ws.Range("I86").Value = 2192.8572
$ws.Range("J86").Value = 5586.5
$ws.Range("K86").Value = 2192.8572
$ws.Range("L86").Value = 5586.5
$ws.Range("M86").Value = -1069.8572
$ws.Range("N86").Value = -7832.5
$ws.Range("H89").Value = 3426.9092
$ws.Range("I89").Value = 2192.8572
$ws.Range("J89").Value = 5586.5
$ws.Range("K89").Value = 10964.286
$ws.Range("L89").Value = 27932.5
$ws.Range("M89").Value = -5348.286
$ws.Range("N89").Value = -39164.5
$ws.Range("H94").Value = 1630
$ws.Range("J94").Value = 1169.6666
$ws.Range("L94").Value = 1169.6666
$ws.Range("N94").Value = -2071.6666
$ws.Range("H107").Value = 752.19446
$ws.Range("I107").Value = 575.1739
$ws.Range("K107").Value = 575.1739
$ws.Range("M107").Value = 1344.8261
$ws.Range("H136").Value = 11620.02
$ws.Range("I136").Value = 6825.0835
$ws.Range("J136").Value = 13095.385
$ws.Range("K136").Value = 20475.2505
$ws.Range("L136").Value = 39286.155
$ws.Range("M136").Value = -17925.2505
$ws.Range("N136").Value = -44386.155

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 1610.95
$ws.Range("I103").Value = 727.5
$ws.Range("K103").Value = 2182.5
$ws.Range("M103").Value = -1303.5
$ws.Range("H126").Value = 9187.25
$ws.Range("I126").Value = 5699.8
$ws.Range("J126").Value = 14999.667
$ws.Range("K126").Value = 17099.4
$ws.Range("L126").Value = 44999.001
$ws.Range("M126").Value = -12159.4
$ws.Range("N126").Value = -54879.001
$ws.Range("H131").Value = 1407.56
$ws.Range("I131").Value = 662.63635
$ws.Range("J131").Value = 1499.6293
$ws.Range("K131").Value = 1987.90905
$ws.Range("L131").Value = 4498.8879
$ws.Range("M131").Value = 3052.09095
$ws.Range("N131").Value = -14578.8879

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1077.7778
$ws.Range("I70").Value = 775.125
$ws.Range("K70").Value = 775.125
$ws.Range("M70").Value = -505.125
$ws.Range("H73").Value = 1077.7778
$ws.Range("I73").Value = 775.125
$ws.Range("K73").Value = 775.125
$ws.Range("M73").Value = 160.875
$ws.Range("H80").Value = 9604.947
$ws.Range("I80").Value = 5884
$ws.Range("J80").Value = 17667
$ws.Range("K80").Value = 5884
$ws.Range("L80").Value = 17667
$ws.Range("M80").Value = -4886
$ws.Range("N80").Value = -19663
$ws.Range("H83").Value = 9604.947
$ws.Range("I83").Value = 5884
$ws.Range("J83").Value = 17667
$ws.Range("K83").Value = 29420
$ws.Range("L83").Value = 88335
$ws.Range("M83").Value = -24428
$ws.Range("N83").Value = -98319
$ws.Range("H97").Value = 1225.5
$ws.Range("I97").Value = 1139.4
$ws.Range("K97").Value = 1139.4
$ws.Range("M97").Value = -643.4000000000001
$ws.Range("H118").Value = 90155
$ws.Range("J118").Value = 90155
$ws.Range("L118").Value = 90155
$ws.Range("N118").Value = -93469
$ws.Range("H132").Value = 16496.889
$ws.Range("I132").Value = 6057.25
$ws.Range("J132").Value = 100014
$ws.Range("K132").Value = 18171.75
$ws.Range("L132").Value = 300042
$ws.Range("M132").Value = -15641.75
$ws.Range("N132").Value = -305102

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 9657.5
$ws.Range("J68").Value = 6896.25
$ws.Range("L68").Value = 6896.25
$ws.Range("N68").Value = -8394.25
$ws.Range("H71").Value = 9657.5
$ws.Range("J71").Value = 6896.25
$ws.Range("L71").Value = 34481.25
$ws.Range("N71").Value = -41969.25
$ws.Range("H93").Value = 7300.5293
$ws.Range("I93").Value = 7353.5835
$ws.Range("J93").Value = 7173.2
$ws.Range("K93").Value = 7353.5835
$ws.Range("L93").Value = 7173.2
$ws.Range("M93").Value = -6105.5835
$ws.Range("N93").Value = -9669.200000000001
$ws.Range("H95").Value = 36999.668
$ws.Range("J95").Value = 32500
$ws.Range("L95").Value = 32500
$ws.Range("N95").Value = -37992
$ws.Range("H100").Value = 2933.8572
$ws.Range("I100").Value = 2546.4211
$ws.Range("J100").Value = 3751.7778
$ws.Range("K100").Value = 2546.4211
$ws.Range("L100").Value = 3751.7778
$ws.Range("M100").Value = -2005.4211
$ws.Range("N100").Value = -4833.7778
$ws.Range("H106").Value = 35500
$ws.Range("J106").Value = 35500
$ws.Range("L106").Value = 35500
$ws.Range("N106").Value = -38024
$ws.Range("H132").Value = 1220212.1
$ws.Range("I132").Value = 4322.923
$ws.Range("J132").Value = 2010540.2
$ws.Range("K132").Value = 12968.769
$ws.Range("L132").Value = 6031620.6
$ws.Range("M132").Value = -10438.769
$ws.Range("N132").Value = -6036680.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2599.7693
$ws.Range("I96").Value = 2032.6666
$ws.Range("J96").Value = 2769.9
$ws.Range("K96").Value = 2032.6666
$ws.Range("L96").Value = 2769.9
$ws.Range("M96").Value = -659.6666
$ws.Range("N96").Value = -5515.9
$ws.Range("H117").Value = 92803
$ws.Range("J117").Value = 92803
$ws.Range("L117").Value = 92803
$ws.Range("N117").Value = -101981
